$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DistTableDisplay6b")

# Update title text: year range changed from 1995-96 to 2005-06
$ws.Range("A3").Value = "AGRICULTURAL CENSUS , 2005-06"

# Update data values in rows 9-24 (columns C, E, F - D stays 0)
$data = @{
    9  = @{ C = 51;  E = 11;   F = 11 }
    10 = @{ C = 141; E = 96;   F = 96 }
    11 = @{ C = 192; E = 107;  F = 107 }
    12 = @{ C = 347; E = 426;  F = 426 }
    13 = @{ C = 347; E = 426;  F = 426 }
    14 = @{ C = 239; E = 352;  F = 352 }
    15 = @{ C = 82;  E = 159;  F = 159 }
    16 = @{ C = 321; E = 511;  F = 511 }
    17 = @{ C = 42;  E = 99;   F = 99 }
    18 = @{ C = 15;  E = 47;   F = 47 }
    19 = @{ C = 6;   E = 15;   F = 15 }
    20 = @{ C = 63;  E = 162;  F = 162 }
    21 = @{ C = 0;   E = 0;    F = 0 }
    22 = @{ C = 1;   E = 56;   F = 56 }
    23 = @{ C = 1;   E = 56;   F = 56 }
    24 = @{ C = 924; E = 1261; F = 1261 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
}
